# ---------------------------------------------------------------------------
# Applies the edits described by the commit diff:
#   1. Remove the stray <w:bookmarkStart/bookmarkEnd w:name="_GoBack"> pair
#      that sits at the very top of the document (before "Министерство...").
#   2. "добавление, ... пользователя;" gains ", уведомления" before the
#      closing semicolon, split across 4 runs.
#   3. "изменение статуса сообщений …;" becomes "изменение статуса
#      привычек;", split across 3 runs, with the _GoBack bookmark now
#      anchored right after "привычек" (before the final ";").
#   4. ", камеры, файлов пользователя." drops "камеры, " -> ", файлов
#      пользователя.", split across 3 runs (leaving the preceding
#      "просмотр Календаря" run untouched).
#   5. The two runs "Изделие, входящее ... предназначенное " / "для
#      помощи ... жизни" are merged back into a single run.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

function Split-AtOffset($startOffset, $endOffset) {
    # Forces a run boundary at $startOffset (and at $endOffset) without
    # altering visible formatting: toggling a character property on and
    # back off causes the engine to break the run(s) spanning the range,
    # and the empty-but-present <w:rPr> collapses back to its original
    # contents once the property is restored.
    $r = $d.Range($startOffset, $endOffset)
    $r.Font.Bold = $true
    $r.Font.Bold = $false
}

# ---------------------------------------------------------------------------
# 1. Remove the original "_GoBack" bookmark near the top of the document.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. "добавление, удаление и редактирование информации о Привычках
#    пользователя;" -> add ", уведомления" before the ";", split into
#    4 runs:
#      "добавление, удаление и редактирование информации о Привычках
#       пользователя" / ", уведомлени" / "я" / ";"
# ---------------------------------------------------------------------------
$rng = $d.Content
$oldText2 = "добавление, удаление и редактирование информации о Привычках пользователя;"
if ($rng.Find.Execute($oldText2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $start = $rng.Start
    $part1 = "добавление, удаление и редактирование информации о Привычках пользователя"
    $part2 = ", уведомлени"
    $part3 = "я"
    $part4 = ";"
    $newText2 = $part1 + $part2 + $part3 + $part4
    $rng.Text = $newText2
    $end = $start + $newText2.Length

    Split-AtOffset ($start + $part1.Length) $end
    Split-AtOffset ($start + $part1.Length + $part2.Length) $end
    Split-AtOffset ($start + $part1.Length + $part2.Length + $part3.Length) $end
}

# ---------------------------------------------------------------------------
# 3. "изменение статуса сообщений …;" -> "изменение статуса привычек;",
#    split into 3 runs: "изменение статуса " / "привычек" / ";" with the
#    _GoBack bookmark inserted right between "привычек" and ";".
# ---------------------------------------------------------------------------
$rng = $d.Content
$oldText3 = "изменение статуса сообщений …;"
if ($rng.Find.Execute($oldText3, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $start = $rng.Start
    $part1 = "изменение статуса "
    $part2 = "привычек"
    $part3 = ";"
    $newText3 = $part1 + $part2 + $part3
    $rng.Text = $newText3
    $end = $start + $newText3.Length

    Split-AtOffset ($start + $part1.Length) $end

    $bookmarkPoint = $d.Range($start + $part1.Length + $part2.Length, $start + $part1.Length + $part2.Length)
    $d.Bookmarks.Add("_GoBack", $bookmarkPoint)
}

# ---------------------------------------------------------------------------
# 4. ", камеры, файлов пользователя." -> ", файлов пользователя.",
#    split into 3 runs: ", " / "файлов" / " пользователя." while the
#    preceding "просмотр Календаря" run stays untouched.
# ---------------------------------------------------------------------------
$rng = $d.Content
$oldText4 = ", камеры, файлов пользователя."
if ($rng.Find.Execute($oldText4, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $start = $rng.Start
    $part1 = ", "
    $part2 = "файлов"
    $part3 = " пользователя."
    $newText4 = $part1 + $part2 + $part3
    $rng.Text = $newText4
    $end = $start + $newText4.Length

    # Break away from the preceding "просмотр Календаря" run first.
    Split-AtOffset $start $end
    Split-AtOffset ($start + $part1.Length) $end
    Split-AtOffset ($start + $part1.Length + $part2.Length) $end
}

# ---------------------------------------------------------------------------
# 5. Merge "Изделие, входящее ... предназначенное " and "для помощи ...
#    жизни" back into a single run (content unchanged).
# ---------------------------------------------------------------------------
$rng = $d.Content
$fullText5 = "Изделие, входящее в курсовой проект – мобильное приложение, предназначенное для помощи пользователям по выработке привычек и улучшения их качества жизни"
if ($rng.Find.Execute($fullText5, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $start = $rng.Start
    $rng.Text = "x"
    $tmpEnd = $start + 1
    $rng2 = $d.Range($start, $tmpEnd)
    $rng2.Text = $fullText5
}

Write-Output "edits applied"
